# data cleanup continued in player_per_game_df
# Remove the "Malcolm Brogdon" row (row 27) entirely from the smoy pivot
# table and bump Leandro Barbosa's (row 25) award count from 1 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sanity-check we are deleting the expected row before mutating anything.
$playerToRemove = $ws.Range("A27").Value()
if ($playerToRemove -ne "Malcolm Brogdon") {
    throw "Expected row 27 to be 'Malcolm Brogdon' but found '$playerToRemove'"
}

# Drop the entire row; rows below shift up by one, and the now-unused
# shared string gets pruned from sharedStrings.xml on save.
$ws.Rows(27).EntireRow.Delete()

# Leandro Barbosa's award count increases from 1 to 2.
$playerToUpdate = $ws.Range("A25").Value()
if ($playerToUpdate -ne "Leandro Barbosa") {
    throw "Expected row 25 to be 'Leandro Barbosa' but found '$playerToUpdate'"
}
$ws.Range("B25").Value = 2
